$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates scraped by the GitHub Actions cron job: refreshed price (D) and
# 1h volume-change (E) for every coin, plus a name/link swap for rows 45-46 (Frax
# and Decentraland traded ranking positions this run). "ForceText" marks D values
# that look like plain numbers (e.g. "1.008", "0.3800") - Coinranking prices are
# scraped as text, so those must be kept as text instead of being auto-converted
# to numbers (which would silently drop significant trailing zeros).
$updates = @(
    @{ Row = 2; D = '27.442.92'; ForceText = $false; E = '  +3.63%  ' },
    @{ Row = 3; D = '1.799.96'; ForceText = $false; E = '  +4.60%  ' },
    @{ Row = 4; D = '1.008'; ForceText = $true; E = '  +0.51%  ' },
    @{ Row = 5; D = '335.62'; ForceText = $true; E = '  +1.12%  ' },
    @{ Row = 6; D = '1.004'; ForceText = $true; E = '  +0.24%  ' },
    @{ Row = 7; D = '0.3800'; ForceText = $true; E = '  +2.51%  ' },
    @{ Row = 8; D = '0.3480'; ForceText = $true; E = '  +3.77%  ' },
    @{ Row = 9; D = '48.77'; ForceText = $true; E = '  +1.29%  ' },
    @{ Row = 10; E = '  +2.51%  ' },
    @{ Row = 11; D = '0.07576'; ForceText = $true; E = '  +2.60%  ' },
    @{ Row = 12; D = '1.005'; ForceText = $true; E = '  +0.32%  ' },
    @{ Row = 13; D = '22.09'; ForceText = $true; E = '  +10.20%  ' },
    @{ Row = 14; D = '6.546'; ForceText = $true; E = '  +2.62%  ' },
    @{ Row = 15; D = '1.800.03'; ForceText = $false; E = '  +4.79%  ' },
    @{ Row = 16; D = '7.100'; ForceText = $true; E = '  +1.29%  ' },
    @{ Row = 17; D = '0.00001103'; ForceText = $true; E = '  +3.29%  ' },
    @{ Row = 18; D = '0.06684'; ForceText = $true; E = '  +1.11%  ' },
    @{ Row = 19; D = '85.06'; ForceText = $true; E = '  +3.57%  ' },
    @{ Row = 20; D = '1.004'; ForceText = $true; E = '  +0.27%  ' },
    @{ Row = 21; D = '17.52'; ForceText = $true; E = '  +6.05%  ' },
    @{ Row = 22; D = '6.489'; ForceText = $true; E = '  +6.25%  ' },
    @{ Row = 23; D = '27.451.34'; ForceText = $false; E = '  +3.86%  ' },
    @{ Row = 24; E = '  -1.42%  ' },
    @{ Row = 25; D = '2.458'; ForceText = $true; E = '  +1.05%  ' },
    @{ Row = 26; D = '2.576'; ForceText = $true; E = '  +7.77%  ' },
    @{ Row = 27; D = '21.57'; ForceText = $true; E = '  +11.43%  ' },
    @{ Row = 28; E = '  +4.98%  ' },
    @{ Row = 29; D = '150.67'; ForceText = $true; E = '  -0.86%  ' },
    @{ Row = 30; D = '2.003.14'; ForceText = $false; E = '  +4.78%  ' },
    @{ Row = 31; D = '134.10'; ForceText = $true; E = '  +2.67%  ' },
    @{ Row = 32; D = '4.080'; ForceText = $true; E = '  -0.93%  ' },
    @{ Row = 33; D = '6.128'; ForceText = $true; E = '  +2.96%  ' },
    @{ Row = 34; D = '0.08660'; ForceText = $true; E = '  +0.63%  ' },
    @{ Row = 35; D = '13.36'; ForceText = $true; E = '  +5.81%  ' },
    @{ Row = 36; D = '1.669'; ForceText = $true; E = '  -1.52%  ' },
    @{ Row = 37; D = '5.519'; ForceText = $true; E = '  +3.48%  ' },
    @{ Row = 38; D = '0.6888'; ForceText = $true; E = '  +11.70%  ' },
    @{ Row = 39; D = '0.2222'; ForceText = $true; E = '  +3.30%  ' },
    @{ Row = 40; D = '0.02372'; ForceText = $true; E = '  +2.32%  ' },
    @{ Row = 41; D = '0.06385'; ForceText = $true; E = '  +3.36%  ' },
    @{ Row = 42; D = '8.878'; ForceText = $true; E = '  +5.34%  ' },
    @{ Row = 43; E = '  +4.36%  ' },
    @{ Row = 44; D = '14.37'; ForceText = $true; E = '  +2.94%  ' },
    @{ Row = 45; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.6464'; ForceText = $true; E = '  +8.79%  ' },
    @{ Row = 46; B = 'Frax'; C = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D = '1.004'; ForceText = $true; E = '  +0.29%  ' },
    @{ Row = 47; D = '3.847'; ForceText = $true; E = '  -1.45%  ' },
    @{ Row = 48; E = '  +5.53%  ' },
    @{ Row = 49; D = '130.82'; ForceText = $true; E = '  +2.22%  ' },
    @{ Row = 50; D = '0.07237'; ForceText = $true; E = '  +1.04%  ' },
    @{ Row = 51; D = '79.85'; ForceText = $true; E = '  +4.34%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$row").Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        if ($u.ForceText) {
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.NumberFormat = "General"
        } else {
            $cell.Value = $u.D
        }
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$row").Value = $u.E }
}